$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.513.05'
$ws.Range('E2').Value = '  +2.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.870.26'
$ws.Range('E3').Value = '  +1.42%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.013'
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.24'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4790'
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3781'
$ws.Range('E8').Value = '  +3.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07361'
$ws.Range('E9').Value = '  +2.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9379'
$ws.Range('E10').Value = '  +1.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.73'
$ws.Range('E11').Value = '  +5.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07852'
$ws.Range('E12').Value = '  +2.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.860.70'
$ws.Range('E13').Value = '  -2.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.444'
$ws.Range('E14').Value = '  +2.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.574'
$ws.Range('E15').Value = '  +2.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.75'
$ws.Range('E16').Value = '  +2.14%  '
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008916'
$ws.Range('E18').Value = '  +3.33%  '
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('E20').Value = '  +2.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.518.73'
$ws.Range('E21').Value = '  +2.09%  '
$ws.Range('E22').Value = '  +1.53%  '
$ws.Range('E23').Value = '  +0.84%  '
$ws.Range('E24').Value = '  +1.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.84'
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('E26').Value = '  +2.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.029'
$ws.Range('E27').Value = '  +1.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '115.88'
$ws.Range('E28').Value = '  +1.51%  '
$ws.Range('E29').Value = '  +1.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08928'
$ws.Range('E30').Value = '  +0.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.334'
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.216'
$ws.Range('E32').Value = '  +3.63%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7568'
$ws.Range('E33').Value = '  +1.49%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.612'
$ws.Range('E34').Value = '  +2.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.714'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  +4.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.119'
$ws.Range('E37').Value = '  +0.95%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05275'
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.999'
$ws.Range('E39').Value = '  +0.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5347'
$ws.Range('E40').Value = '  +2.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.092'
$ws.Range('E41').Value = '  +1.95%  '
$ws.Range('E42').Value = '  +1.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.502'
$ws.Range('E43').Value = '  +3.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.72'
$ws.Range('E44').Value = '  +1.80%  '
$ws.Range('E45').Value = '  +1.97%  '
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.662'
$ws.Range('E47').Value = '  +3.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.97'
$ws.Range('E48').Value = '  +1.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '67.48'
$ws.Range('E49').Value = '  +2.04%  '
$ws.Range('E50').Value = '  +1.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9247'
$ws.Range('E51').Value = '  +4.60%  '
